$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.186.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "'3.135.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.22%  "
$ws.Range("D5").Value = "'524.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.44%  "
$ws.Range("D6").Value = "'134.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.90%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'3.131.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.30%  "
$ws.Range("E9").Value = "  -4.53%  "
$ws.Range("D10").Value = "'7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.36%  "
$ws.Range("E11").Value = "  -8.31%  "
$ws.Range("E12").Value = "  -6.47%  "
$ws.Range("D13").Value = "'3.670.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "'25.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("D16").Value = "'3.131.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.18%  "
$ws.Range("D17").Value = "'58.166.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("E18").Value = "  -6.40%  "
$ws.Range("E19").Value = "  -5.75%  "
$ws.Range("D20").Value = "'13.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("D21").Value = "'7.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").Value = "'343.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.59%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("D25").Value = "'67.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.91%  "
$ws.Range("D26").Value = "'3.261.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").Value = "'0.172"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'0.0₃0951"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.20%  "
$ws.Range("D29").Value = "'0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Value = "'6.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -7.24%  "
$ws.Range("D33").Value = "'6.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.21%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").Value = "'157.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("E39").Value = "  -10.27%  "
$ws.Range("E40").Value = "  -5.28%  "
$ws.Range("D41").Value = "'3.165.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("D42").Value = "'40.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").Value = "'23.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.09%  "
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("D46").Value = "'3.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'2.283.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -7.41%  "
$ws.Range("D50").Value = "'6.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("E51").Value = "  -1.93%  "
